$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9 (Leve Item ID 5487)
$ws.Range("H9").Value = 7059.6665
$ws.Range("I9").Value = 8576
$ws.Range("J9").Value = 994.3333
$ws.Range("K9").Value = 8576
$ws.Range("L9").Value = 994.3333
$ws.Range("M9").Value = -8407
$ws.Range("N9").Value = -1332.3333

# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 364.35
$ws.Range("I33").Value = 364
$ws.Range("K33").Value = 364
$ws.Range("M33").Value = -135

# Row 92 (Leve Item ID 19901)
$ws.Range("H92").Value = 1022
$ws.Range("I92").Value = 817.55554
$ws.Range("K92").Value = 817.55554
$ws.Range("M92").Value = 430.44446

# Row 136 (Leve Item ID 42164)
$ws.Range("H136").Value = 88329.336
$ws.Range("J136").Value = 88329.336
$ws.Range("L136").Value = 88329.336
$ws.Range("N136").Value = -98529.336

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 5804.3335
$ws.Range("J137").Value = 1981.3334
$ws.Range("L137").Value = 5944.0002
$ws.Range("N137").Value = -11044.0002


$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 15636135
$ws.Range("I32").Value = 23260894
$ws.Range("K32").Value = 23260894
$ws.Range("M32").Value = -23260607

# Row 41 (Leve Item ID 2501)
$ws.Range("H41").Value = 6079.1816
$ws.Range("I41").Value = 6079.1816
$ws.Range("K41").Value = 6079.1816
$ws.Range("M41").Value = -5665.1816

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 43482340
$ws.Range("I61").Value = 76925544
$ws.Range("K61").Value = 76925544
$ws.Range("M61").Value = -76925332

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 41717230
$ws.Range("I74").Value = 71510780
$ws.Range("J74").Value = 6258.6
$ws.Range("K74").Value = 71510780
$ws.Range("L74").Value = 6258.6
$ws.Range("M74").Value = -71509906
$ws.Range("N74").Value = -8006.6

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 41717230
$ws.Range("I77").Value = 71510780
$ws.Range("J77").Value = 6258.6
$ws.Range("K77").Value = 357553900
$ws.Range("L77").Value = 31293
$ws.Range("M77").Value = -357549532
$ws.Range("N77").Value = -40029

# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 850.3929000000001
$ws.Range("I97").Value = 351.69565
$ws.Range("K97").Value = 351.69565
$ws.Range("M97").Value = 144.30435

# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 2188.5
$ws.Range("I102").Value = 2188.5
$ws.Range("K102").Value = 2188.5
$ws.Range("M102").Value = -566.5

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 2426.9429
$ws.Range("I122").Value = 1353
$ws.Range("K122").Value = 4059
$ws.Range("M122").Value = -1609

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 62503580
$ws.Range("I132").Value = 3816.9333
$ws.Range("K132").Value = 11450.7999
$ws.Range("M132").Value = -8920.7999

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 43482340
$ws.Range("I136").Value = 76925544
$ws.Range("K136").Value = 230776632
$ws.Range("M136").Value = -230774082


$ws = $wb.Worksheets.Item("BSM")
# Row 22 (Leve Item ID 5092)
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

# Row 80 (Leve Item ID 13747)
$ws.Range("H80").Value = 2741.111
$ws.Range("I80").Value = 898.6667
$ws.Range("J80").Value = 4583.5557
$ws.Range("K80").Value = 898.6667
$ws.Range("L80").Value = 4583.5557
$ws.Range("M80").Value = 99.33330000000001
$ws.Range("N80").Value = -6579.5557

# Row 83 (Leve Item ID 13747)
$ws.Range("H83").Value = 2741.111
$ws.Range("I83").Value = 898.6667
$ws.Range("J83").Value = 4583.5557
$ws.Range("K83").Value = 4493.3335
$ws.Range("L83").Value = 22917.7785
$ws.Range("M83").Value = 498.6665000000003
$ws.Range("N83").Value = -32901.7785

# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 1624.7142
$ws.Range("I94").Value = 1474.8
$ws.Range("K94").Value = 1474.8
$ws.Range("M94").Value = -1023.8

# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 1633.8334
$ws.Range("I105").Value = 1464.875
$ws.Range("J105").Value = 1971.75
$ws.Range("K105").Value = 1464.875
$ws.Range("L105").Value = 1971.75
$ws.Range("M105").Value = 282.125
$ws.Range("N105").Value = -5465.75

# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 6148.5264
$ws.Range("I107").Value = 4433.3
$ws.Range("K107").Value = 4433.3
$ws.Range("M107").Value = -2513.3

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 2744
$ws.Range("I134").Value = 2642.9033
$ws.Range("K134").Value = 7928.7099
$ws.Range("M134").Value = -5393.7099


$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 43484036
$ws.Range("I31").Value = 5153.4614
$ws.Range("K31").Value = 5153.4614
$ws.Range("M31").Value = -4858.4614

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 43484036
$ws.Range("I34").Value = 5153.4614
$ws.Range("K34").Value = 5153.4614
$ws.Range("M34").Value = -4951.4614

# Row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 3118.7368
$ws.Range("I86").Value = 2596.3845
$ws.Range("K86").Value = 2596.3845
$ws.Range("M86").Value = -1473.3845

# Row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 3118.7368
$ws.Range("I89").Value = 2596.3845
$ws.Range("K89").Value = 12981.9225
$ws.Range("M89").Value = -7365.922500000001

# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 9792.223
$ws.Range("I99").Value = 10448.286
$ws.Range("K99").Value = 10448.286
$ws.Range("M99").Value = -8950.286

# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 9216
$ws.Range("J105").Value = 27752.5
$ws.Range("L105").Value = 27752.5
$ws.Range("N105").Value = -31246.5

# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 1591
$ws.Range("I122").Value = 1523.3077
$ws.Range("K122").Value = 4569.9231
$ws.Range("M122").Value = -2119.9231

# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 9792.223
$ws.Range("I126").Value = 10448.286
$ws.Range("K126").Value = 31344.858
$ws.Range("M126").Value = -28874.858


$ws = $wb.Worksheets.Item("CUL")
# Row 12 (Leve Item ID 4854)
$ws.Range("H12").Value = 506.42856
$ws.Range("I12").Value = 304.25
$ws.Range("J12").Value = 587.3
$ws.Range("K12").Value = 912.75
$ws.Range("L12").Value = 1761.9
$ws.Range("M12").Value = -739.75
$ws.Range("N12").Value = -2107.9

# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 1640.1177
$ws.Range("J122").Value = 1940.6666
$ws.Range("L122").Value = 17465.9994
$ws.Range("N122").Value = -22365.9994

# Row 137 (Leve Item ID 44088)
$ws.Range("H137").Value = 3998.8823
$ws.Range("J137").Value = 5810.778
$ws.Range("L137").Value = 17432.334
$ws.Range("N137").Value = -27632.334


$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 4230.8
$ws.Range("I80").Value = 3925
$ws.Range("J80").Value = 4361.857
$ws.Range("K80").Value = 3925
$ws.Range("L80").Value = 4361.857
$ws.Range("M80").Value = -2927
$ws.Range("N80").Value = -6357.857

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 4230.8
$ws.Range("I83").Value = 3925
$ws.Range("J83").Value = 4361.857
$ws.Range("K83").Value = 19625
$ws.Range("L83").Value = 21809.285
$ws.Range("M83").Value = -14633
$ws.Range("N83").Value = -31793.285

# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 2790
$ws.Range("I122").Value = 1983.3334
$ws.Range("K122").Value = 5950.0002
$ws.Range("M122").Value = -3500.0002

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 1831.5454
$ws.Range("I132").Value = 1266.5385
$ws.Range("J132").Value = 2647.6667
$ws.Range("K132").Value = 3799.6155
$ws.Range("L132").Value = 7943.000100000001
$ws.Range("M132").Value = -1269.6155
$ws.Range("N132").Value = -13003.0001


$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 418.13043
$ws.Range("I16").Value = 459
$ws.Range("J16").Value = 145.66667
$ws.Range("K16").Value = 459
$ws.Range("L16").Value = 145.66667
$ws.Range("M16").Value = -289
$ws.Range("N16").Value = -485.66667

# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2853.1177
$ws.Range("I22").Value = 1291.7778
$ws.Range("K22").Value = 1291.7778
$ws.Range("M22").Value = -996.7778000000001

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2853.1177
$ws.Range("I27").Value = 1291.7778
$ws.Range("K27").Value = 1291.7778
$ws.Range("M27").Value = -1184.7778

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 2281
$ws.Range("J46").Value = 3823.6
$ws.Range("L46").Value = 3823.6
$ws.Range("N46").Value = -4199.6

# Row 47 (Leve Item ID 3138)
$ws.Range("H47").Value = 24780
$ws.Range("J47").Value = 24780
$ws.Range("L47").Value = 24780
$ws.Range("N47").Value = -25760

# Row 52 (Leve Item ID 3138)
$ws.Range("H52").Value = 24780
$ws.Range("J52").Value = 24780
$ws.Range("L52").Value = 24780
$ws.Range("N52").Value = -25246

# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 3769.2856
$ws.Range("I61").Value = 2125
$ws.Range("J61").Value = 5961.6665
$ws.Range("K61").Value = 2125
$ws.Range("L61").Value = 5961.6665
$ws.Range("M61").Value = -1923
$ws.Range("N61").Value = -6365.6665

# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 2729.6365
$ws.Range("J100").Value = 3076.923
$ws.Range("L100").Value = 3076.923
$ws.Range("N100").Value = -4158.923

# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 3769.2856
$ws.Range("I113").Value = 2125
$ws.Range("J113").Value = 5961.6665
$ws.Range("K113").Value = 2125
$ws.Range("L113").Value = 5961.6665
$ws.Range("M113").Value = 45
$ws.Range("N113").Value = -10301.6665


$ws = $wb.Worksheets.Item("WVR")
# Row 96 (Leve Item ID 19977)
$ws.Range("H96").Value = 4777.7144
$ws.Range("J96").Value = 5273.2
$ws.Range("L96").Value = 5273.2
$ws.Range("N96").Value = -8019.2

